# Apply "update publikace + difuzni-prezkum" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("publikace")

# New shared strings must be introduced in this order to reproduce the
# target uniqueCount ordering: 28-43, sbornik70, sbornik71, 31-40,
# m-468.html, dm-468.pdf

# Row 17: Kdy se senator musi stat poslancem? - add Pages (I17) and Fulltext (N17)
$ws.Range("I17").Value = "28-43"
$ws.Range("N17").Value = "https://www.wintr.cz/images/sborniky/sbornik70.pdf"

# Row 16: Kdy je proporcionalni branit sebeposkozovani? - add Pages (I16) and Fulltext (N16)
$ws.Range("N16").Value = "https://www.wintr.cz/images/sborniky/sbornik71.pdf"
$ws.Range("I16").Value = "31-40"

# Row 9: Usneseni vlady jako pravni predpis... (Jurisprudence, 33-38)
# add Url (H9) and Fulltext (N9)
$ws.Range("H9").Value = "https://www.jurisprudence.cz/cz/casopis/usneseni-vlady-jako-pravni-predpis-terminologicky-zmatek-sui-generis.m-468.html"
$ws.Range("N9").Value = "https://www.jurisprudence.cz/cz/casopis/usneseni-vlady-jako-pravni-predpis-terminologicky-zmatek-sui-generis.dm-468.pdf"

# Update view/selection to match saved state
$ws.Range("N9").Select()
